# Update the cryptos price-list sheet: latest Price (D) / Volume(1h) (E)
# snapshot, including the EnergySwap/BabyDogeCoin row-order swap at rows 49-50.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store these as text (matching the
# source file's inline-string cells) instead of auto-converting numeric-
# looking values like '23.06' or '1.000' into floats and losing formatting.
$q = "'"

$ws.Range("D2").Value = $q + '29.051.48'
$ws.Range("E2").Value = '  +0.47%  '

$ws.Range("D3").Value = $q + '1.832.49'
$ws.Range("E3").Value = '  +0.53%  '

$ws.Range("D4").Value = $q + '0.9986'
$ws.Range("E4").Value = '  +0.46%  '

$ws.Range("D5").Value = $q + '241.74'
$ws.Range("E5").Value = '  -0.63%  '

$ws.Range("D6").Value = $q + '0.6188'
$ws.Range("E6").Value = '  -1.58%  '

$ws.Range("D7").Value = $q + '1.000'
$ws.Range("E7").Value = '  +0.36%  '

$ws.Range("D8").Value = $q + '0.07449'
$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").Value = $q + '0.2927'
$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("D10").Value = $q + '23.06'

$ws.Range("D11").Value = $q + '0.07668'
$ws.Range("E11").Value = '  -0.05%  '

$ws.Range("D12").Value = $q + '1.845.32'
$ws.Range("E12").Value = '  +1.08%  '

$ws.Range("D13").Value = $q + '5.004'
$ws.Range("E13").Value = '  +0.59%  '

$ws.Range("D14").Value = $q + '0.6734'
$ws.Range("E14").Value = '  +1.18%  '

$ws.Range("D15").Value = $q + '82.90'
$ws.Range("E15").Value = '  +0.06%  '

$ws.Range("D16").Value = $q + '0.000009150'
$ws.Range("E16").Value = '  -5.66%  '

$ws.Range("D17").Value = $q + '5.912'
$ws.Range("E17").Value = '  -1.71%  '

$ws.Range("D18").Value = $q + '29.059.78'
$ws.Range("E18").Value = '  +0.36%  '

$ws.Range("D19").Value = $q + '2.086.09'
$ws.Range("E19").Value = '  +0.85%  '

$ws.Range("D20").Value = $q + '240.64'
$ws.Range("E20").Value = '  +7.00%  '

$ws.Range("D21").Value = $q + '12.71'
$ws.Range("E21").Value = '  +1.52%  '

$ws.Range("D22").Value = $q + '1.001'
$ws.Range("E22").Value = '  +0.63%  '

$ws.Range("D23").Value = $q + '7.213'
$ws.Range("E23").Value = '  +1.48%  '

$ws.Range("D24").Value = $q + '1.0000'
$ws.Range("E24").Value = '  +0.43%  '

$ws.Range("D25").Value = $q + '158.72'
$ws.Range("E25").Value = '  -0.86%  '

$ws.Range("D26").Value = $q + '0.1413'
$ws.Range("E26").Value = '  +0.48%  '

$ws.Range("D27").Value = $q + '8.503'
$ws.Range("E27").Value = '  +0.24%  '

$ws.Range("D28").Value = $q + '17.89'
$ws.Range("E28").Value = '  +0.27%  '

$ws.Range("D29").Value = $q + '1.500'
$ws.Range("E29").Value = '  +0.38%  '

$ws.Range("D30").Value = $q + '0.05624'
$ws.Range("E30").Value = '  +3.28%  '

$ws.Range("D31").Value = $q + '4.115'
$ws.Range("E31").Value = '  +1.76%  '

$ws.Range("D32").Value = $q + '4.131'
$ws.Range("E32").Value = '  +0.60%  '

$ws.Range("D33").Value = $q + '1.201'
$ws.Range("E33").Value = '  +0.40%  '

$ws.Range("D34").Value = $q + '1.842'
$ws.Range("E34").Value = '  -0.26%  '

$ws.Range("D35").Value = $q + '0.7412'
$ws.Range("E35").Value = '  -0.10%  '

$ws.Range("D36").Value = $q + '1.143'
$ws.Range("E36").Value = '  +1.02%  '

$ws.Range("D37").Value = $q + '2.656'
$ws.Range("E37").Value = '  +1.97%  '

$ws.Range("D38").Value = $q + '2.771'
$ws.Range("E38").Value = '  +1.07%  '

$ws.Range("D39").Value = $q + '0.01785'
$ws.Range("E39").Value = '  +0.55%  '

$ws.Range("D40").Value = $q + '1.211.17'
$ws.Range("E40").Value = '  -2.27%  '

$ws.Range("D41").Value = $q + '6.395'
$ws.Range("E41").Value = '  -3.83%  '

$ws.Range("D42").Value = $q + '0.8987'
$ws.Range("E42").Value = '  +0.27%  '

$ws.Range("D43").Value = $q + '0.9987'
$ws.Range("E43").Value = '  +0.39%  '

$ws.Range("D44").Value = $q + '101.40'
$ws.Range("E44").Value = '  +0.18%  '

$ws.Range("D45").Value = $q + '1.983.95'
$ws.Range("E45").Value = '  +0.72%  '

$ws.Range("D46").Value = $q + '65.43'
$ws.Range("E46").Value = '  +0.97%  '

$ws.Range("D47").Value = $q + '0.5086'
$ws.Range("E47").Value = '  +0.45%  '

$ws.Range("D48").Value = $q + '0.4061'
$ws.Range("E48").Value = '  +0.36%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = $q + '0.00000000118'
$ws.Range("E49").Value = '  -3.58%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = $q + '9.143'
$ws.Range("E50").Value = '  +2.43%  '

$ws.Range("D51").Value = $q + '0.05813'
$ws.Range("E51").Value = '  +0.64%  '
